$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Extend the row-number column (A) down through row 37, copying the
#    existing header-cell formatting (bold/border/center) from A19.
$ws.Range("A19").Copy()
$ws.Range("A20:A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false
for ($r = 20; $r -le 37; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# 2) Refresh the screener ticker lists in columns B-F (rows 2-37) to the
#    newly uploaded data. Only cells whose content actually changed are
#    touched; unchanged placeholder cells are left as-is.
$ws.Range("B2").Value = "NSE:APEX"
$ws.Range("C2").Value = "NSE:AAREYDRUGS"
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "NSE:AARTIIND"
$ws.Range("F2").Value = "NSE:INFY"
$ws.Range("B3").Value = "NSE:HDFCGROWTH"
$ws.Range("C3").Value = "NSE:ADSL"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "NSE:ADANIENT"
$ws.Range("F3").Value = "NSE:ITC"
$ws.Range("B4").Value = "NSE:HDFCNIFTY"
$ws.Range("C4").Value = "NSE:ATAM"
$ws.Range("E4").Value = "NSE:ASHOKLEY"
$ws.Range("F4").ClearContents()
$ws.Range("B5").Value = "NSE:INFY"
$ws.Range("E5").Value = "NSE:AUROPHARMA"
$ws.Range("F5").ClearContents()
$ws.Range("B6").Value = "NSE:ITC"
$ws.Range("C6").Value = "NSE:BSL"
$ws.Range("E6").Value = "NSE:CANBK"
$ws.Range("B7").Value = "NSE:KARURVYSYA"
$ws.Range("C7").Value = "NSE:CAREERP"
$ws.Range("E7").Value = "NSE:CONCOR"
$ws.Range("B8").Value = "NSE:KREBSBIO"
$ws.Range("C8").Value = "NSE:CONCOR"
$ws.Range("E8").Value = "NSE:GRANULES"
$ws.Range("B9").Value = "NSE:ORCHPHARMA"
$ws.Range("C9").Value = "NSE:DBCORP"
$ws.Range("E9").Value = "NSE:HINDALCO"
$ws.Range("B10").Value = "NSE:POWERMECH"
$ws.Range("C10").Value = "NSE:DEN"
$ws.Range("E10").Value = "NSE:IDFC"
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = "NSE:DPSCLTD"
$ws.Range("E11").Value = "NSE:IRCTC"
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = "NSE:ENERGYDEV"
$ws.Range("E12").Value = "NSE:JINDALSTEL"
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "NSE:EUROTEXIND"
$ws.Range("E13").Value = "NSE:LAURUSLABS"
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "NSE:GEOJITFSL"
$ws.Range("E14").Value = "NSE:NATIONALUM"
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = "NSE:GODFRYPHLP"
$ws.Range("E15").Value = "NSE:NTPC"
$ws.Range("C16").Value = "NSE:GREENLAM"
$ws.Range("E16").Value = "NSE:PFC"
$ws.Range("C17").Value = "NSE:GVKPIL"
$ws.Range("E17").Value = "NSE:PNB"
$ws.Range("C18").Value = "NSE:HAVISHA"
$ws.Range("E18").Value = "NSE:POLYCAB"
$ws.Range("C19").Value = "NSE:HLVLTD"
$ws.Range("E19").Value = "NSE:POWERGRID"
$ws.Range("C20").Value = "NSE:IFCI"
$ws.Range("E20").Value = "NSE:RECLTD"
$ws.Range("C21").Value = "NSE:IRMENERGY"
$ws.Range("C22").Value = "NSE:KELLTONTEC"
$ws.Range("C23").Value = "NSE:LATENTVIEW"
$ws.Range("C24").Value = "NSE:M&MFIN"
$ws.Range("C25").Value = "NSE:MADRASFERT"
$ws.Range("C26").Value = "NSE:MAHABANK"
$ws.Range("C27").Value = "NSE:MINDTECK"
$ws.Range("C28").Value = "NSE:MIRZAINT"
$ws.Range("C29").Value = "NSE:MMFL"
$ws.Range("C30").Value = "NSE:MONTECARLO"
$ws.Range("C31").Value = "NSE:NCLIND"
$ws.Range("C32").Value = "NSE:NFL"
$ws.Range("C33").Value = "NSE:NILAINFRA"
$ws.Range("C34").Value = "NSE:NLCINDIA"
$ws.Range("C35").Value = "NSE:PRAJIND"
$ws.Range("C36").Value = "NSE:RBL"
$ws.Range("C37").Value = "NSE:RUCHIRA"
